$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark ------------------------------------
#        It used to sit right after the picture later in the document;
#        it now belongs right after "IDG PA28" in the title, i.e. right
#        before the trailing "X" that is about to be dropped.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Paragraphs(1).Range.Duplicate
$bmRange.Start = $bmRange.Start + 8
$bmRange.End = $bmRange.Start
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 2. Rename "IDG PA28X" -> "IDG PA28" --------------------------------
#        Delete just the trailing "X" (now right after the bookmark) so
#        the run split made by the bookmark, and the untouched rsid on
#        the following run, are both preserved.
$xRange = $d.Paragraphs(1).Range.Duplicate
$xRange.Start = $d.Bookmarks("_GoBack").Start
$xRange.End = $xRange.Start + 1
$xRange.Delete()

# --- 3. Merge the two runs that spell out how to open the panel --------
#        "The panel can be shown by selecting " + "<ldquo>Mini Panel<rdquo> from
#        the Utilities menu." become one run with identical text.
$d.Content.Find.Execute("Utilities menu.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Utilities menu.", 2) | Out-Null

# --- 4. Merge "Show Aircraft Conf" + "iguration" into one run ----------
$d.Content.Find.Execute("Show Aircraft Conf", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Show Aircraft Conf", 2) | Out-Null
